# feat: add 2022-Q1 data
#
# The workbook previously ended with a single "总计" (summary) sheet.
# This script:
#   1. Duplicates the "总计" sheet (so both copies start out with the same
#      formatting/sheet properties as the original).
#   2. Repurposes the first copy (in place of the old "总计") as the new
#      "2022-Q1" sheet, replacing its contents with the per-fund holding
#      details for that quarter.
#   3. Repurposes the second copy as the new "总计" sheet, updating the
#      summary table with a new 2022-Q1 row at the top (and every other
#      row pushed down one).
#
# Final tab order: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计

$wb = $excel.ActiveWorkbook

$old = $wb.Worksheets.Item("总计")

# Duplicate "总计" -> creates "总计 (2)" immediately after it, carrying over
# the same sheetPr / column-A and header styling.
$old.Copy($null, $old)
$newTotal = $wb.Worksheets.Item("总计 (2)")

# ------------------------------------------------------------------
# Turn the original "总计" sheet into "2022-Q1" (fund holding details).
# ------------------------------------------------------------------
$old.Name = "2022-Q1"
$old.Cells.Clear()

$srcHdr = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$srcHdr.Copy()
$old.Range("B1:H1").PasteSpecial(-4122)

$srcA = $wb.Worksheets.Item("2021-Q4").Range("A2:A5")
$srcA.Copy()
$old.Range("A2:A5").PasteSpecial(-4122)

$old.Range("B1").Value = "基金代码"
$old.Range("C1").Value = "基金名称"
$old.Range("D1").Value = "基金规模"
$old.Range("E1").Value = "股票总仓位"
$old.Range("F1").Value = "仓位占比"
$old.Range("G1").Value = "持有市值(亿元)"
$old.Range("H1").Value = "仓位排名"

$old.Range("A2").Value = 0
$old.Range("A3").Value = 1
$old.Range("A4").Value = 2
$old.Range("A5").Value = 3

# Columns B-G hold text (fund codes / formatted numbers-as-text), force
# the cells to text so values like "002666" keep their leading zero and
# "11.96" is not re-interpreted as a number.
$old.Range("B2:G5").NumberFormat = "@"

$old.Range("B2").Value = "002666"
$old.Range("C2").Value = "前海开源沪港深创新成长灵活配置混合A"
$old.Range("D2").Value = "11.96"
$old.Range("E2").Value = "81.64"
$old.Range("F2").Value = "5.79"
$old.Range("G2").Value = "0.6925"
$old.Range("H2").Value = 10

$old.Range("B3").Value = "002667"
$old.Range("C3").Value = "前海开源沪港深创新成长灵活配置混合C"
$old.Range("D3").Value = "3.25"
$old.Range("E3").Value = "81.64"
$old.Range("F3").Value = "5.79"
$old.Range("G3").Value = "0.1882"
$old.Range("H3").Value = 10

$old.Range("B4").Value = "004403"
$old.Range("C4").Value = "平安股息精选沪港深股票A"
$old.Range("D4").Value = "1.16"
$old.Range("E4").Value = "93.51"
$old.Range("F4").Value = "3.95"
$old.Range("G4").Value = "0.0458"
$old.Range("H4").Value = 10

$old.Range("B5").Value = "004404"
$old.Range("C5").Value = "平安股息精选沪港深股票C"
$old.Range("D5").Value = "0.01"
$old.Range("E5").Value = "93.51"
$old.Range("F5").Value = "3.95"
$old.Range("G5").Value = "0.0004"
$old.Range("H5").Value = 10

# ------------------------------------------------------------------
# Turn the duplicated sheet into the updated "总计" summary, with a new
# 2022-Q1 row inserted at the top (row 2) and everything else shifted
# down by one row.
# ------------------------------------------------------------------

# Extend the index-column / row styling down to the new row 7 (copy the
# formatting already present on row 6 of the duplicate).
$newTotal.Range("A6").Copy()
$newTotal.Range("A7").PasteSpecial(-4122)

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 4
$newTotal.Range("D2").Value = 0.93

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 5
$newTotal.Range("D3").Value = 1.45

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 7
$newTotal.Range("D4").Value = 1.58

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 15
$newTotal.Range("D5").Value = 22

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 20
$newTotal.Range("D6").Value = 34.74

$newTotal.Range("A7").Value = 5
$newTotal.Range("B7").Value = "2020-Q4"
$newTotal.Range("C7").Value = 22
$newTotal.Range("D7").Value = 21.32

$newTotal.Name = "总计"
